# Update "想去人数" (F column) values across the four worksheets to reflect
# the latest generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 247
$ws1.Range("F4").Value  = 839
$ws1.Range("F5").Value  = 242
$ws1.Range("F6").Value  = 413
$ws1.Range("F7").Value  = 594
$ws1.Range("F10").Value = 346
$ws1.Range("F11").Value = 146
$ws1.Range("F12").Value = 663
$ws1.Range("F13").Value = 87
$ws1.Range("F14").Value = 1812
$ws1.Range("F15").Value = 355
$ws1.Range("F16").Value = 3283
$ws1.Range("F17").Value = 325
$ws1.Range("F18").Value = 494
$ws1.Range("F19").Value = 55
$ws1.Range("F20").Value = 143
$ws1.Range("F21").Value = 131

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 3

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5336
$ws3.Range("F3").Value = 324
$ws3.Range("F4").Value = 265

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5336
$ws4.Range("F4").Value  = 324
$ws4.Range("F6").Value  = 265
$ws4.Range("F7").Value  = 247
$ws4.Range("F14").Value = 839
$ws4.Range("F17").Value = 242
$ws4.Range("F18").Value = 413
$ws4.Range("F19").Value = 594
$ws4.Range("F23").Value = 346
$ws4.Range("F24").Value = 146
$ws4.Range("F27").Value = 663
$ws4.Range("F28").Value = 87
$ws4.Range("F30").Value = 1812
$ws4.Range("F31").Value = 355
$ws4.Range("F32").Value = 3284
$ws4.Range("F34").Value = 325
$ws4.Range("F35").Value = 494
$ws4.Range("F36").Value = 55
$ws4.Range("F37").Value = 3
$ws4.Range("F38").Value = 143
$ws4.Range("F40").Value = 131
